# Auto-generated: apply scheduled market-data refresh to Sheets/Ultros_Profits.xlsx
# Updates plain numeric values (currentAveragePrice*, LevePrice*, LeveProfit*) per-row
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets; no formulas/styles are touched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1282.871
$ws.Range("I15").Value = 1282.871
$ws.Range("K15").Value = 3848.613
$ws.Range("M15").Value = -3679.613
$ws.Range("H17").Value = 2079
$ws.Range("J17").Value = 2079
$ws.Range("L17").Value = 6237
$ws.Range("N17").Value = -6573
$ws.Range("H74").Value = 6422.154
$ws.Range("I74").Value = 3355.4285
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 3355.4285
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -2419.4285
$ws.Range("N74").Value = -11872
$ws.Range("H76").Value = 4999.385
$ws.Range("I76").Value = 5320.3
$ws.Range("J76").Value = 3929.6667
$ws.Range("K76").Value = 5320.3
$ws.Range("L76").Value = 3929.6667
$ws.Range("M76").Value = -5005.3
$ws.Range("N76").Value = -4559.6667
$ws.Range("H77").Value = 6422.154
$ws.Range("I77").Value = 3355.4285
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 16777.1425
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -12097.1425
$ws.Range("N77").Value = -59360
$ws.Range("H79").Value = 4999.385
$ws.Range("I79").Value = 5320.3
$ws.Range("J79").Value = 3929.6667
$ws.Range("K79").Value = 5320.3
$ws.Range("L79").Value = 3929.6667
$ws.Range("M79").Value = -4228.3
$ws.Range("N79").Value = -6113.6667
$ws.Range("H132").Value = 15458.6
$ws.Range("I132").Value = 1202.9
$ws.Range("K132").Value = 3608.7
$ws.Range("M132").Value = -1078.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14462
$ws.Range("I2").Value = 23633.928
$ws.Range("J2").Value = 1621.3
$ws.Range("K2").Value = 23633.928
$ws.Range("L2").Value = 1621.3
$ws.Range("M2").Value = -23520.928
$ws.Range("N2").Value = -1847.3
$ws.Range("H44").Value = 40000
$ws.Range("J44").Value = 40000
$ws.Range("L44").Value = 40000
$ws.Range("N44").Value = -40976
$ws.Range("H62").Value = 55000
$ws.Range("J62").Value = 55000
$ws.Range("L62").Value = 55000
$ws.Range("N62").Value = -56248
$ws.Range("H63").Value = 7334.8335
$ws.Range("I63").Value = 8502.25
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 8502.25
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -7816.25
$ws.Range("N63").Value = -6372
$ws.Range("H65").Value = 55000
$ws.Range("J65").Value = 55000
$ws.Range("L65").Value = 165000
$ws.Range("N65").Value = -171240
$ws.Range("H66").Value = 7334.8335
$ws.Range("I66").Value = 8502.25
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 42511.25
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -39079.25
$ws.Range("N66").Value = -31864
$ws.Range("H97").Value = 1612.5
$ws.Range("I97").Value = 1850
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 1850
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -1354
$ws.Range("N97").Value = -1892
$ws.Range("H116").Value = 14462
$ws.Range("I116").Value = 23633.928
$ws.Range("J116").Value = 1621.3
$ws.Range("K116").Value = 23633.928
$ws.Range("L116").Value = 1621.3
$ws.Range("M116").Value = -21339.928
$ws.Range("N116").Value = -6209.3
$ws.Range("H132").Value = 4336.6665
$ws.Range("I132").Value = 3930.762
$ws.Range("J132").Value = 5757.3335
$ws.Range("K132").Value = 11792.286
$ws.Range("L132").Value = 17272.0005
$ws.Range("M132").Value = -9262.286
$ws.Range("N132").Value = -22332.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14462
$ws.Range("I3").Value = 23633.928
$ws.Range("J3").Value = 1621.3
$ws.Range("K3").Value = 23633.928
$ws.Range("L3").Value = 1621.3
$ws.Range("M3").Value = -23519.928
$ws.Range("N3").Value = -1849.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1846.7059
$ws.Range("I58").Value = 1253.5454
$ws.Range("J58").Value = 2934.1667
$ws.Range("K58").Value = 1253.5454
$ws.Range("L58").Value = 2934.1667
$ws.Range("M58").Value = -1050.5454
$ws.Range("N58").Value = -3340.1667
$ws.Range("H59").Value = 115000
$ws.Range("J59").Value = 58000
$ws.Range("L59").Value = 58000
$ws.Range("N59").Value = -60290
$ws.Range("H112").Value = 39812.5
$ws.Range("J112").Value = 39812.5
$ws.Range("L112").Value = 39812.5
$ws.Range("N112").Value = -42766.5
$ws.Range("H136").Value = 1846.7059
$ws.Range("I136").Value = 1253.5454
$ws.Range("J136").Value = 2934.1667
$ws.Range("K136").Value = 3760.6362
$ws.Range("L136").Value = 8802.500100000001
$ws.Range("M136").Value = -1210.6362
$ws.Range("N136").Value = -13902.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3840.6875
$ws.Range("J39").Value = 4070.0667
$ws.Range("L39").Value = 12210.2001
$ws.Range("N39").Value = -12798.2001
$ws.Range("H140").Value = 50001850
$ws.Range("I140").Value = 62501564
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 187504692
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = -187499512
$ws.Range("N140").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 61143.633
$ws.Range("I80").Value = 81648.07000000001
$ws.Range("J80").Value = 3731.2
$ws.Range("K80").Value = 81648.07000000001
$ws.Range("L80").Value = 3731.2
$ws.Range("M80").Value = -80650.07000000001
$ws.Range("N80").Value = -5727.2
$ws.Range("H83").Value = 61143.633
$ws.Range("I83").Value = 81648.07000000001
$ws.Range("J83").Value = 3731.2
$ws.Range("K83").Value = 408240.35
$ws.Range("L83").Value = 18656
$ws.Range("M83").Value = -403248.35
$ws.Range("N83").Value = -28640
$ws.Range("H132").Value = 5516.524
$ws.Range("J132").Value = 5246.5
$ws.Range("L132").Value = 15739.5
$ws.Range("N132").Value = -20799.5
$ws.Range("H138").Value = 72000
$ws.Range("J138").Value = 72000
$ws.Range("L138").Value = 72000
$ws.Range("N138").Value = -82280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3258.3333
$ws.Range("I61").Value = 1387.5
$ws.Range("K61").Value = 1387.5
$ws.Range("M61").Value = -1185.5
$ws.Range("H100").Value = 104569.37
$ws.Range("I100").Value = 140495.38
$ws.Range("K100").Value = 140495.38
$ws.Range("M100").Value = -139954.38
$ws.Range("H113").Value = 3258.3333
$ws.Range("I113").Value = 1387.5
$ws.Range("K113").Value = 1387.5
$ws.Range("M113").Value = 782.5
$ws.Range("H132").Value = 3613.7673
$ws.Range("I132").Value = 2881.1724
$ws.Range("K132").Value = 8643.5172
$ws.Range("M132").Value = -6113.5172

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2669.6667
$ws.Range("I132").Value = 2410.7878
$ws.Range("K132").Value = 7232.3634
$ws.Range("M132").Value = -4702.3634
$ws.Range("H139").Value = 60000
$ws.Range("J139").Value = 60000
$ws.Range("L139").Value = 60000
$ws.Range("N139").Value = -70280
$ws.Range("H140").Value = 99999
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360

